# Chức năng quản lý tài khoản
# Adds a new member row (row 29) to the members sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 29

$ws.Cells.Item($row, 1).Value = 29            # A29 - Key
$ws.Cells.Item($row, 3).Value = "default.jpg" # C29 - AVT
$ws.Cells.Item($row, 8).Value = "N/A"         # H29 - Chuyen nganh
$ws.Cells.Item($row, 9).Value = "N/A"         # I29 - Dai hoc
$ws.Cells.Item($row, 10).Value = "N/A"        # J29 - SDT
$ws.Cells.Item($row, 11).Value = "N/A"        # K29 - Mail
$ws.Cells.Item($row, 12).Value = "N/A"        # L29 - Dia chi
$ws.Cells.Item($row, 13).Value = "Chưa có"    # M29 - Don vi
$ws.Cells.Item($row, 14).Value = "Chưa có"    # N29 - Chuc vu
$ws.Cells.Item($row, 15).Value = $false       # O29 - La LT?
$ws.Cells.Item($row, 16).Value = $false       # P29 - Qua PTBT?
